$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text so they remain
# stored the same way (as text) as the rest of the Price/Volume columns.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D15", "D17", "D20", "D21", "D22", "D25", "D28", "D30", "D32", "D35", "D39", "D40", "D41", "D42", "D45", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Update the Price (D) and Volume(1h) (E) columns with the latest snapshot values.
$ws.Range("D2").Value = "26.154.22"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.612.52"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "213.14"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "0.482"
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "0.0619"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").Value = "18.39"
$ws.Range("E10").Value = "  +2.94%  "
$ws.Range("D11").Value = "0.0797"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "1.838.08"
$ws.Range("E12").Value = "  +0.94%  "
$ws.Range("D13").Value = "1.612.11"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "26.160.06"
$ws.Range("D17").Value = "60.87"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  +1.91%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "198.53"
$ws.Range("E20").Value = "  +4.80%  "
$ws.Range("D21").Value = "4.26"
$ws.Range("E21").Value = "  +2.13%  "
$ws.Range("D22").Value = "9.48"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "142.41"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").Value = "15.21"
$ws.Range("E28").Value = "  +2.12%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").Value = "3.15"
$ws.Range("E32").Value = "  +2.22%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("D35").Value = "2.34"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("D36").Value = "1.107.77"
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "0.506"
$ws.Range("E39").Value = "  +2.71%  "
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "0.793"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "0.801"
$ws.Range("E42").Value = "  +8.14%  "
$ws.Range("D43").Value = "1.749.89"
$ws.Range("E43").Value = "  +0.99%  "
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").Value = "93.15"
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("E46").Value = "  +6.76%  "
$ws.Range("E47").Value = "  +8.55%  "
$ws.Range("D48").Value = "53.93"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E50").Value = "  -0.24%  "
